$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 1398.3334
$ws.Cells.Item(2, 9).Value = 182.66667
$ws.Cells.Item(2, 10).Value = 3829.6667
$ws.Cells.Item(2, 11).Value = 182.66667
$ws.Cells.Item(2, 12).Value = 3829.6667
$ws.Cells.Item(2, 13).Value = -69.66667000000001
$ws.Cells.Item(2, 14).Value = -4055.6667
$ws.Cells.Item(17, 8).Value = 1492.3182
$ws.Cells.Item(17, 10).Value = 1492.3182
$ws.Cells.Item(17, 12).Value = 4476.9546
$ws.Cells.Item(17, 14).Value = -4812.9546
$ws.Cells.Item(18, 8).Value = 370
$ws.Cells.Item(18, 9).Value = 370
$ws.Cells.Item(18, 11).Value = 370
$ws.Cells.Item(18, 13).Value = -86
$ws.Cells.Item(40, 8).Value = 4954.773
$ws.Cells.Item(40, 9).Value = 2398.5
$ws.Cells.Item(40, 10).Value = 6415.5
$ws.Cells.Item(40, 11).Value = 2398.5
$ws.Cells.Item(40, 12).Value = 6415.5
$ws.Cells.Item(40, 13).Value = -2223.5
$ws.Cells.Item(40, 14).Value = -6765.5
$ws.Cells.Item(46, 8).Value = 4522
$ws.Cells.Item(46, 10).Value = 4440
$ws.Cells.Item(46, 12).Value = 13320
$ws.Cells.Item(46, 14).Value = -13558
$ws.Cells.Item(60, 8).Value = 4522
$ws.Cells.Item(60, 10).Value = 4440
$ws.Cells.Item(60, 12).Value = 13320
$ws.Cells.Item(60, 14).Value = -14288
$ws.Cells.Item(88, 8).Value = 4285.4287
$ws.Cells.Item(88, 9).Value = 4999.5
$ws.Cells.Item(88, 10).Value = 3333.3333
$ws.Cells.Item(88, 11).Value = 4999.5
$ws.Cells.Item(88, 12).Value = 3333.3333
$ws.Cells.Item(88, 13).Value = -4593.5
$ws.Cells.Item(88, 14).Value = -4145.3333
$ws.Cells.Item(91, 8).Value = 4285.4287
$ws.Cells.Item(91, 9).Value = 4999.5
$ws.Cells.Item(91, 10).Value = 3333.3333
$ws.Cells.Item(91, 11).Value = 4999.5
$ws.Cells.Item(91, 12).Value = 3333.3333
$ws.Cells.Item(91, 13).Value = -3595.5
$ws.Cells.Item(91, 14).Value = -6141.3333
$ws.Cells.Item(98, 8).Value = 1165.7142
$ws.Cells.Item(98, 9).Value = 920.2
$ws.Cells.Item(98, 11).Value = 920.2
$ws.Cells.Item(98, 13).Value = 577.8
$ws.Cells.Item(122, 8).Value = 1165.7142
$ws.Cells.Item(122, 9).Value = 920.2
$ws.Cells.Item(122, 11).Value = 2760.6
$ws.Cells.Item(122, 13).Value = -310.6000000000004
$ws.Cells.Item(128, 8).Value = 90000
$ws.Cells.Item(128, 10).Value = 90000
$ws.Cells.Item(128, 12).Value = 90000
$ws.Cells.Item(128, 14).Value = -99960
$ws.Cells.Item(130, 8).Value = 110000
$ws.Cells.Item(130, 10).Value = 110000
$ws.Cells.Item(130, 12).Value = 110000
$ws.Cells.Item(130, 14).Value = -120040
$ws.Cells.Item(132, 8).Value = 4035.087
$ws.Cells.Item(132, 9).Value = 3943.1904
$ws.Cells.Item(132, 11).Value = 11829.5712
$ws.Cells.Item(132, 13).Value = -9299.5712
$ws.Cells.Item(137, 8).Value = 3291.3704
$ws.Cells.Item(137, 9).Value = 1975.1111
$ws.Cells.Item(137, 11).Value = 5925.3333
$ws.Cells.Item(137, 13).Value = -3375.3333
$ws.Cells.Item(138, 8).Value = 5404.4136
$ws.Cells.Item(138, 10).Value = 6982.4287
$ws.Cells.Item(138, 12).Value = 20947.2861
$ws.Cells.Item(138, 14).Value = -31227.2861

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4927.9043
$ws.Cells.Item(32, 10).Value = 30862.25
$ws.Cells.Item(32, 12).Value = 30862.25
$ws.Cells.Item(32, 14).Value = -31436.25
$ws.Cells.Item(61, 8).Value = 4691.8096
$ws.Cells.Item(61, 9).Value = 1738
$ws.Cells.Item(61, 11).Value = 1738
$ws.Cells.Item(61, 13).Value = -1526
$ws.Cells.Item(74, 8).Value = 5244.4707
$ws.Cells.Item(74, 10).Value = 10682.2
$ws.Cells.Item(74, 12).Value = 10682.2
$ws.Cells.Item(74, 14).Value = -12430.2
$ws.Cells.Item(77, 8).Value = 5244.4707
$ws.Cells.Item(77, 10).Value = 10682.2
$ws.Cells.Item(77, 12).Value = 53411
$ws.Cells.Item(77, 14).Value = -62147
$ws.Cells.Item(97, 8).Value = 441.2353
$ws.Cells.Item(97, 9).Value = 352.64285
$ws.Cells.Item(97, 11).Value = 352.64285
$ws.Cells.Item(97, 13).Value = 143.35715
$ws.Cells.Item(136, 8).Value = 4691.8096
$ws.Cells.Item(136, 9).Value = 1738
$ws.Cells.Item(136, 11).Value = 5214
$ws.Cells.Item(136, 13).Value = -2664

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 7872.88
$ws.Cells.Item(31, 9).Value = 4361.857
$ws.Cells.Item(31, 10).Value = 9238.277
$ws.Cells.Item(31, 11).Value = 4361.857
$ws.Cells.Item(31, 12).Value = 9238.277
$ws.Cells.Item(31, 13).Value = -4066.857
$ws.Cells.Item(31, 14).Value = -9828.277
$ws.Cells.Item(34, 8).Value = 7872.88
$ws.Cells.Item(34, 9).Value = 4361.857
$ws.Cells.Item(34, 10).Value = 9238.277
$ws.Cells.Item(34, 11).Value = 4361.857
$ws.Cells.Item(34, 12).Value = 9238.277
$ws.Cells.Item(34, 13).Value = -4159.857
$ws.Cells.Item(34, 14).Value = -9642.277
$ws.Cells.Item(52, 8).Value = 98765
$ws.Cells.Item(52, 10).Value = 98765
$ws.Cells.Item(52, 12).Value = 98765
$ws.Cells.Item(52, 14).Value = -99353
$ws.Cells.Item(86, 8).Value = 6124.6665
$ws.Cells.Item(86, 9).Value = 4933
$ws.Cells.Item(86, 10).Value = 7316.3335
$ws.Cells.Item(86, 11).Value = 4933
$ws.Cells.Item(86, 12).Value = 7316.3335
$ws.Cells.Item(86, 13).Value = -3810
$ws.Cells.Item(86, 14).Value = -9562.333500000001
$ws.Cells.Item(89, 8).Value = 6124.6665
$ws.Cells.Item(89, 9).Value = 4933
$ws.Cells.Item(89, 10).Value = 7316.3335
$ws.Cells.Item(89, 11).Value = 24665
$ws.Cells.Item(89, 12).Value = 36581.6675
$ws.Cells.Item(89, 13).Value = -19049
$ws.Cells.Item(89, 14).Value = -47813.6675
$ws.Cells.Item(132, 8).Value = 3729.4443
$ws.Cells.Item(132, 9).Value = 2890
$ws.Cells.Item(132, 11).Value = 8670
$ws.Cells.Item(132, 13).Value = -6140

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 708.7143
$ws.Cells.Item(5, 9).Value = 598.6667
$ws.Cells.Item(5, 10).Value = 791.25
$ws.Cells.Item(5, 11).Value = 1796.0001
$ws.Cells.Item(5, 12).Value = 2373.75
$ws.Cells.Item(5, 13).Value = -1684.0001
$ws.Cells.Item(5, 14).Value = -2597.75
$ws.Cells.Item(12, 8).Value = 387.84616
$ws.Cells.Item(12, 9).Value = 113
$ws.Cells.Item(12, 10).Value = 470.3
$ws.Cells.Item(12, 11).Value = 339
$ws.Cells.Item(12, 12).Value = 1410.9
$ws.Cells.Item(12, 13).Value = -166
$ws.Cells.Item(12, 14).Value = -1756.9
$ws.Cells.Item(113, 8).Value = 1212.7142
$ws.Cells.Item(113, 9).Value = 388.5
$ws.Cells.Item(113, 10).Value = 1542.4
$ws.Cells.Item(113, 11).Value = 1165.5
$ws.Cells.Item(113, 12).Value = 4627.200000000001
$ws.Cells.Item(113, 13).Value = 1004.5
$ws.Cells.Item(113, 14).Value = -8967.200000000001
$ws.Cells.Item(135, 8).Value = 708.7143
$ws.Cells.Item(135, 9).Value = 598.6667
$ws.Cells.Item(135, 10).Value = 791.25
$ws.Cells.Item(135, 11).Value = 5388.0003
$ws.Cells.Item(135, 12).Value = 7121.25
$ws.Cells.Item(135, 13).Value = -2853.0003
$ws.Cells.Item(135, 14).Value = -12191.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 92832.5
$ws.Cells.Item(7, 9).Value = 131999.25
$ws.Cells.Item(7, 10).Value = 14499
$ws.Cells.Item(7, 11).Value = 131999.25
$ws.Cells.Item(7, 12).Value = 14499
$ws.Cells.Item(7, 13).Value = -131887.25
$ws.Cells.Item(7, 14).Value = -14723
$ws.Cells.Item(40, 8).Value = 4332
$ws.Cells.Item(40, 10).Value = 7766
$ws.Cells.Item(40, 12).Value = 7766
$ws.Cells.Item(40, 14).Value = -8038
$ws.Cells.Item(46, 8).Value = 3278.318
$ws.Cells.Item(46, 10).Value = 3694.9333
$ws.Cells.Item(46, 12).Value = 3694.9333
$ws.Cells.Item(46, 14).Value = -4070.9333
$ws.Cells.Item(68, 8).Value = 6198.7334
$ws.Cells.Item(68, 9).Value = 7535.091
$ws.Cells.Item(68, 10).Value = 2523.75
$ws.Cells.Item(68, 11).Value = 7535.091
$ws.Cells.Item(68, 12).Value = 2523.75
$ws.Cells.Item(68, 13).Value = -6786.091
$ws.Cells.Item(68, 14).Value = -4021.75
$ws.Cells.Item(71, 8).Value = 6198.7334
$ws.Cells.Item(71, 9).Value = 7535.091
$ws.Cells.Item(71, 10).Value = 2523.75
$ws.Cells.Item(71, 11).Value = 37675.455
$ws.Cells.Item(71, 12).Value = 12618.75
$ws.Cells.Item(71, 13).Value = -33931.455
$ws.Cells.Item(71, 14).Value = -20106.75
$ws.Cells.Item(100, 8).Value = 7782.375
$ws.Cells.Item(100, 9).Value = 7751.2856
$ws.Cells.Item(100, 10).Value = 8000
$ws.Cells.Item(100, 11).Value = 7751.2856
$ws.Cells.Item(100, 12).Value = 8000
$ws.Cells.Item(100, 13).Value = -7210.2856
$ws.Cells.Item(100, 14).Value = -9082
$ws.Cells.Item(122, 8).Value = 7030.9165
$ws.Cells.Item(122, 9).Value = 6397.364
$ws.Cells.Item(122, 11).Value = 19192.092
$ws.Cells.Item(122, 13).Value = -16742.092
$ws.Cells.Item(126, 8).Value = 92832.5
$ws.Cells.Item(126, 9).Value = 131999.25
$ws.Cells.Item(126, 10).Value = 14499
$ws.Cells.Item(126, 11).Value = 395997.75
$ws.Cells.Item(126, 12).Value = 43497
$ws.Cells.Item(126, 13).Value = -393527.75
$ws.Cells.Item(126, 14).Value = -48437
$ws.Cells.Item(132, 8).Value = 6660.4414
$ws.Cells.Item(132, 9).Value = 5908.857
$ws.Cells.Item(132, 10).Value = 10167.833
$ws.Cells.Item(132, 11).Value = 17726.571
$ws.Cells.Item(132, 12).Value = 30503.499
$ws.Cells.Item(132, 13).Value = -15196.571
$ws.Cells.Item(132, 14).Value = -35563.499

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 55160.89
$ws.Cells.Item(46, 10).Value = 55160.89
$ws.Cells.Item(46, 12).Value = 55160.89
$ws.Cells.Item(46, 14).Value = -55622.89
$ws.Cells.Item(122, 8).Value = 3698.7
$ws.Cells.Item(122, 9).Value = 3726.125
$ws.Cells.Item(122, 11).Value = 11178.375
$ws.Cells.Item(122, 13).Value = -8728.375
$ws.Cells.Item(132, 8).Value = 2781.6
$ws.Cells.Item(132, 9).Value = 1489.1428
$ws.Cells.Item(132, 10).Value = 5797.3335
$ws.Cells.Item(132, 11).Value = 4467.428400000001
$ws.Cells.Item(132, 12).Value = 17392.0005
$ws.Cells.Item(132, 13).Value = -1937.428400000001
$ws.Cells.Item(132, 14).Value = -22452.0005
$ws.Cells.Item(134, 8).Value = 55160.89
$ws.Cells.Item(134, 10).Value = 55160.89
$ws.Cells.Item(134, 12).Value = 165482.67
$ws.Cells.Item(134, 14).Value = -170552.67

Write-Host "Applied Moogle_Profits updates"